$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price/Volume columns for the affected rows so that
# values such as "1.00" or "7.51" are stored as text (matching the original
# inline-string cells) instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '60.776.88'
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").Value = '3.390.53'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '571.03'
$ws.Range("E5").Value = '  -2.15%  '
$ws.Range("E6").Value = '  -3.69%  '
$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").Value = '3.392.26'
$ws.Range("E7").Value = '  -1.92%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").Value = '7.51'
$ws.Range("E10").Value = '  -2.54%  '
$ws.Range("E11").Value = '  -2.10%  '
$ws.Range("D12").Value = '0.394'
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("D13").Value = '3.970.53'
$ws.Range("E13").Value = '  -1.91%  '
$ws.Range("D14").Value = '28.27'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("D17").Value = '3.396.24'
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("D18").Value = '60.885.39'
$ws.Range("E18").Value = '  -1.49%  '
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").Value = '14.03'
$ws.Range("E20").Value = '  -2.49%  '
$ws.Range("D21").Value = '9.03'
$ws.Range("E21").Value = '  -5.73%  '
$ws.Range("D22").Value = '387.64'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '0.561'
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").Value = '73.64'
$ws.Range("E24").Value = '  +0.29%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '0.0000117'
$ws.Range("E26").Value = '  -4.35%  '
$ws.Range("D27").Value = '3.529.41'
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("D28").Value = '0.178'
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  -4.74%  '
$ws.Range("D31").Value = '7.99'
$ws.Range("E31").Value = '  -2.81%  '
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("E33").Value = '  -6.69%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '23.70'
$ws.Range("E35").Value = '  -2.21%  '
$ws.Range("D36").Value = '6.95'
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("D37").Value = '167.27'
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("D38").Value = '3.420.25'
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("D39").Value = '4.97'
$ws.Range("E39").Value = '  -3.77%  '
$ws.Range("E40").Value = '  -5.20%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '27.97'
$ws.Range("E41").Value = '  +3.39%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").Value = '0.0778'
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("D43").Value = '0.782'
$ws.Range("E43").Value = '  -3.42%  '
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = '4.43'
$ws.Range("E45").Value = '  -1.94%  '
$ws.Range("D46").Value = '41.67'
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("D47").Value = '1.68'
$ws.Range("E47").Value = '  -2.78%  '
$ws.Range("D48").Value = '2.556.77'
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("E49").Value = '  -3.92%  '
$ws.Range("D50").Value = '23.31'
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("D51").Value = '6.84'
$ws.Range("E51").Value = '  -1.39%  '
